# Apply the commit: re-order a handful of existing match rows (the scraper
# re-fetched results in a slightly different order for matches sharing the
# same kickoff time) and append 4 newly scraped matches at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Rows 13-15: 3-way rotation of match data (columns F:V).
#    new13 <- old15, new14 <- old13, new15 <- old14
# ---------------------------------------------------------------------
$a13 = $ws.Range("F13:V13").Value()
$a14 = $ws.Range("F14:V14").Value()
$a15 = $ws.Range("F15:V15").Value()

$ws.Range("F13:V13").Value = $a15
$ws.Range("F14:V14").Value = $a13
$ws.Range("F15:V15").Value = $a14

# ---------------------------------------------------------------------
# 2) Rows 19-25: 7-way rotation of match data (columns F:V).
#    new19<-old25, new20<-old19, new21<-old24, new22<-old23,
#    new23<-old22, new24<-old21, new25<-old20
# ---------------------------------------------------------------------
$a19 = $ws.Range("F19:V19").Value()
$a20 = $ws.Range("F20:V20").Value()
$a21 = $ws.Range("F21:V21").Value()
$a22 = $ws.Range("F22:V22").Value()
$a23 = $ws.Range("F23:V23").Value()
$a24 = $ws.Range("F24:V24").Value()
$a25 = $ws.Range("F25:V25").Value()

$ws.Range("F19:V19").Value = $a25
$ws.Range("F20:V20").Value = $a19
$ws.Range("F21:V21").Value = $a24
$ws.Range("F22:V22").Value = $a23
$ws.Range("F23:V23").Value = $a22
$ws.Range("F24:V24").Value = $a21
$ws.Range("F25:V25").Value = $a20

# ---------------------------------------------------------------------
# 3) Rows 60 & 62: swap match data (columns F:V). Row 61 is untouched.
# ---------------------------------------------------------------------
$a60 = $ws.Range("F60:V60").Value()
$a62 = $ws.Range("F62:V62").Value()

$ws.Range("F60:V60").Value = $a62
$ws.Range("F62:V62").Value = $a60

# ---------------------------------------------------------------------
# 4) Append 4 brand-new match rows (124-127), cloning the formatting of
#    the last existing row (123) and then overwriting the values.
# ---------------------------------------------------------------------
$ws.Range("A123:V123").Copy($ws.Range("A124:V124"))
$ws.Range("A123:V123").Copy($ws.Range("A125:V125"))
$ws.Range("A123:V123").Copy($ws.Range("A126:V126"))
$ws.Range("A123:V123").Copy($ws.Range("A127:V127"))

# Row 124: Chrudim 1 x 0 Vlasim
$ws.Range("A124").Value = 123
$ws.Range("E124").Value = 45241.42708333334
$ws.Range("F124").Value = "Chrudim"
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = "Vlasim"
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 2.51
$ws.Range("K124").Value = "09/11/2023 09:13"
$ws.Range("L124").Value = 2.33
$ws.Range("M124").Value = "11/11/2023 10:14"
$ws.Range("N124").Value = 3.3
$ws.Range("O124").Value = "09/11/2023 09:13"
$ws.Range("P124").Value = 3.46
$ws.Range("Q124").Value = "11/11/2023 10:05"
$ws.Range("R124").Value = 2.51
$ws.Range("S124").Value = "09/11/2023 09:13"
$ws.Range("T124").Value = 2.94
$ws.Range("U124").Value = "11/11/2023 10:14"
$ws.Range("V124").Value = "https://www.betexplorer.com/football/czech-republic/fnl/chrudim-vlasim/x8ArAwm6/"

# Row 125: Vyskov 2 x 0 Dukla Prague
$ws.Range("A125").Value = 124
$ws.Range("E125").Value = 45241.42708333334
$ws.Range("F125").Value = "Vyskov"
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = "Dukla Prague"
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2.33
$ws.Range("K125").Value = "09/11/2023 09:13"
$ws.Range("L125").Value = 2.01
$ws.Range("M125").Value = "11/11/2023 09:51"
$ws.Range("N125").Value = 3.6
$ws.Range("O125").Value = "09/11/2023 09:13"
$ws.Range("P125").Value = 3.59
$ws.Range("Q125").Value = "11/11/2023 09:51"
$ws.Range("R125").Value = 2.59
$ws.Range("S125").Value = "09/11/2023 09:13"
$ws.Range("T125").Value = 3.54
$ws.Range("U125").Value = "11/11/2023 09:51"
$ws.Range("V125").Value = "https://www.betexplorer.com/football/czech-republic/fnl/mfk-vyskov-dukla-prague/KWFwBJXa/"

# Row 126: Prostejov 0 x 1 Kromeriz
$ws.Range("A126").Value = 125
$ws.Range("E126").Value = 45241.4375
$ws.Range("F126").Value = "Prostejov"
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = "Kromeriz"
$ws.Range("I126").Value = 1
$ws.Range("J126").Value = 1.93
$ws.Range("K126").Value = "09/11/2023 09:13"
$ws.Range("L126").Value = 1.88
$ws.Range("M126").Value = "11/11/2023 10:20"
$ws.Range("N126").Value = 3.5
$ws.Range("O126").Value = "09/11/2023 09:13"
$ws.Range("P126").Value = 3.55
$ws.Range("Q126").Value = "11/11/2023 10:20"
$ws.Range("R126").Value = 3.34
$ws.Range("S126").Value = "09/11/2023 09:13"
$ws.Range("T126").Value = 4.08
$ws.Range("U126").Value = "11/11/2023 10:20"
$ws.Range("V126").Value = "https://www.betexplorer.com/football/czech-republic/fnl/prostejov-kromeriz/nJHVCu3m/"

# Row 127: Brno 0 x 1 Sigma Olomouc B
$ws.Range("A127").Value = 126
$ws.Range("E127").Value = 45241.625
$ws.Range("F127").Value = "Brno"
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = "Sigma Olomouc B"
$ws.Range("I127").Value = 1
$ws.Range("J127").Value = 1.51
$ws.Range("K127").Value = "09/11/2023 09:13"
$ws.Range("L127").Value = 1.46
$ws.Range("M127").Value = "11/11/2023 14:51"
$ws.Range("N127").Value = 4.12
$ws.Range("O127").Value = "09/11/2023 09:13"
$ws.Range("P127").Value = 4.45
$ws.Range("Q127").Value = "11/11/2023 14:51"
$ws.Range("R127").Value = 5.07
$ws.Range("S127").Value = "09/11/2023 09:13"
$ws.Range("T127").Value = 6.73
$ws.Range("U127").Value = "11/11/2023 14:51"
$ws.Range("V127").Value = "https://www.betexplorer.com/football/czech-republic/fnl/brno-sigma-olomouc/ttHZBaIg/"

Write-Host "done"
